# "Registrar Oferta del fabricante"
# Adds three new error-log entries to the "Errores" sheet, pushing the
# two trailing rows (AuctionManagementBean / BussinessException) down
# so the new rows can be inserted right after the existing list.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Errores")

# Capture the two trailing entries (currently rows 10 and 12) before
# they get overwritten, then relocate them down to rows 13 and 15 to
# make room for the three new rows (9, 10, 11).
$b10 = $ws.Range("B10").Value2
$c10 = $ws.Range("C10").Value2
$b12 = $ws.Range("B12").Value2
$c12 = $ws.Range("C12").Value2

$ws.Range("B10").ClearContents()
$ws.Range("C10").ClearContents()
$ws.Range("B12").ClearContents()
$ws.Range("C12").ClearContents()

$ws.Range("B13").Value2 = $b10
$ws.Range("C13").Value2 = $c10
$ws.Range("B15").Value2 = $b12
$ws.Range("C15").Value2 = $c12

# New rows describing the fabricante-offer registration fixes.
$ws.Range("C9").Value2 = "El unico criterio actual es el mejor precio, el metodo de dar ganador no tiene sentido, el mejor se asigna cada vez que se registra una nueva oferta"
$ws.Range("C10").Value2 = "Al registrar la oferta no se estaba asignando a la subasta correspondiente"
$ws.Range("C11").Value2 = "El metodo de registrar oferta estaba en el bean pero no en el web service y por tanto nunca era empleado"

$ws.Range("B2:C15").Select()
$ws.Range("C12").Select()
